# Auto-generated update of cryptos list (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.293.64"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "1.871.47"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7085"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07773"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.31%  "

# Row 9
$ws.Range("E9").Value = "  -0.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08392"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12
$ws.Range("D12").Value = "1.858.69"
$ws.Range("E12").Value = "  -0.40%  "

# Row 13
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "

# Row 16
$ws.Range("D16").Value = "29.304.95"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.070"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008172"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "

# Row 20
$ws.Range("E20").Value = "  +0.77%  "

# Row 21
$ws.Range("D21").Value = "2.113.84"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.749"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.38%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1591"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.24"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.004"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "

# Row 29
$ws.Range("E29").Value = "  +0.28%  "

# Row 30
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("E31").Value = "  -1.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.300"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05335"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.48%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.939"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "

# Row 35
$ws.Range("E35").Value = "  +0.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.64%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "

# Row 39
$ws.Range("D39").Value = "1.230.52"
$ws.Range("E39").Value = "  +6.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.563"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.16%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8842"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("D46").Value = "2.016.23"
$ws.Range("E46").Value = "  +0.27%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.794"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.81%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5190"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "

# Row 49
$ws.Range("E49").Value = "  +2.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.400"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4310"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
